$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so Excel does not
# auto-coerce these numeric-looking strings into numbers (which would
# drop formatting like trailing zeros, thousands separators, etc.)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '44.909.25'
$ws.Range("E2").Value = '  -4.03%  '

$ws.Range("D3").Value = '2.362.92'
$ws.Range("E3").Value = '  +2.87%  '

$ws.Range("E4").Value = '  +0.36%  '

$ws.Range("D5").Value = '293.34'
$ws.Range("E5").Value = '  -3.26%  '

$ws.Range("D6").Value = '94.48'
$ws.Range("E6").Value = '  -6.76%  '

$ws.Range("D7").Value = '0.556'
$ws.Range("E7").Value = '  -2.08%  '

$ws.Range("E8").Value = '  +0.17%  '

$ws.Range("D9").Value = '0.499'
$ws.Range("E9").Value = '  -4.47%  '

$ws.Range("D10").Value = '34.04'
$ws.Range("E10").Value = '  -6.08%  '

$ws.Range("D11").Value = '0.0774'
$ws.Range("E11").Value = '  -1.84%  '

$ws.Range("D12").Value = '6.93'
$ws.Range("E12").Value = '  -4.82%  '

$ws.Range("E13").Value = '  +0.40%  '

$ws.Range("D14").Value = '2.739.53'
$ws.Range("E14").Value = '  +3.52%  '

$ws.Range("D15").Value = '2.381.26'
$ws.Range("E15").Value = '  +3.63%  '

$ws.Range("D16").Value = '13.87'
$ws.Range("E16").Value = '  +0.60%  '

$ws.Range("D17").Value = '0.819'
$ws.Range("E17").Value = '  +1.38%  '

$ws.Range("D18").Value = '44.899.22'
$ws.Range("E18").Value = '  -3.91%  '

$ws.Range("D19").Value = '12.36'
$ws.Range("E19").Value = '  -4.93%  '

$ws.Range("D20").Value = '0.0₃0926'
$ws.Range("E20").Value = '  -1.19%  '

$ws.Range("D21").Value = '6.05'
$ws.Range("E21").Value = '  +1.08%  '

$ws.Range("D22").Value = '65.96'
$ws.Range("E22").Value = '  +0.26%  '

$ws.Range("D23").Value = '237.50'
$ws.Range("E23").Value = '  -4.95%  '

$ws.Range("D24").Value = '2.73'
$ws.Range("E24").Value = '  -4.93%  '

$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  +0.10%  '

$ws.Range("D26").Value = '1.87'
$ws.Range("E26").Value = '  -1.99%  '

$ws.Range("D27").Value = '2.21'
$ws.Range("E27").Value = '  +0.86%  '

$ws.Range("D28").Value = '37.44'
$ws.Range("E28").Value = '  -10.82%  '

$ws.Range("D29").Value = '9.51'
$ws.Range("E29").Value = '  -3.43%  '

$ws.Range("D30").Value = '3.75'
$ws.Range("E30").Value = '  +15.61%  '

$ws.Range("D31").Value = '20.74'
$ws.Range("E31").Value = '  +3.87%  '

$ws.Range("D32").Value = '148.24'
$ws.Range("E32").Value = '  +0.60%  '

$ws.Range("D33").Value = '2.70'
$ws.Range("E33").Value = '  -4.72%  '

$ws.Range("D34").Value = '5.39'
$ws.Range("E34").Value = '  -3.02%  '

$ws.Range("D35").Value = '0.0756'
$ws.Range("E35").Value = '  -3.84%  '

$ws.Range("E36").Value = '  -3.30%  '

$ws.Range("D37").Value = '1.92'
$ws.Range("E37").Value = '  +9.46%  '

$ws.Range("D38").Value = '0.114'
$ws.Range("E38").Value = '  -2.28%  '

$ws.Range("D39").Value = '14.65'
$ws.Range("E39").Value = '  -8.16%  '

$ws.Range("D40").Value = '3.72'
$ws.Range("E40").Value = '  -5.71%  '

$ws.Range("D41").Value = '0.0292'
$ws.Range("E41").Value = '  -2.26%  '

$ws.Range("D42").Value = '3.15'
$ws.Range("E42").Value = '  -5.21%  '

$ws.Range("D43").Value = '1.925.39'
$ws.Range("E43").Value = '  +6.04%  '

$ws.Range("E44").Value = '  +0.20%  '

$ws.Range("D45").Value = '88.97'
$ws.Range("E45").Value = '  -1.56%  '

$ws.Range("D46").Value = '1.69'

$ws.Range("D47").Value = '8.57'
$ws.Range("E47").Value = '  +9.25%  '

$ws.Range("D48").Value = '15.10'
$ws.Range("E48").Value = '  +15.31%  '

$ws.Range("D49").Value = '99.42'
$ws.Range("E49").Value = '  +4.11%  '

$ws.Range("D50").Value = '2.611.61'
$ws.Range("E50").Value = '  +3.56%  '

$ws.Range("D51").Value = '0.182'
$ws.Range("E51").Value = '  -5.55%  '
